# Update "想去人数" (F column) figures for refreshed scrape data
# in both the "展览" and "全部类型" worksheets.

$wb = $excel.ActiveWorkbook

# Map: worksheet name -> hashtable of row -> new value
$updates = @{
    "展览"     = @{ 2 = 135; 3 = 2164; 5 = 11413; 9 = 11356; 13 = 1745; 14 = 5664; 16 = 3480; 17 = 175 }
    "全部类型" = @{ 2 = 135; 3 = 2164; 7 = 11413; 11 = 11356; 15 = 1745; 17 = 5664; 19 = 3480; 20 = 175 }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Cells.Item([int]$row, 6).Value = $rows[$row]
    }
}
